$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D ("Company") to make room for "Linkedin".
# This shifts: D->E, E->F, F->G, G->H, H->I (and copies D1 header style).
$ws.Columns("D").Insert()

# Update header row
$ws.Range("D1").Value = "Linkedin"
$ws.Range("F1").Value = "Job"
$ws.Range("G1").Value = "Location"
$ws.Range("H1").Value = "Phones"
$ws.Range("I1").Value = "Emails"

# Populate the new Linkedin column and reconcile the Emails column per row
$ws.Range("D2").Value = "https://linkedin.com/in/aidan-mccarron-92416356"
$ws.Range("I2").Value = ""

$ws.Range("D3").Value = "https://linkedin.com/in/john-higgins-mciob-b25b37a3"
$ws.Range("I3").Value = ""

$ws.Range("D4").Value = "https://linkedin.com/in/michael-yohanis-mciob-060b5147"
$ws.Range("I4").Value = ""

$ws.Range("D5").Value = "https://linkedin.com/in/declan-mclogan-cmiosh-ll-m-48482613"
$ws.Range("I5").Value = ""

$ws.Range("D6").Value = "https://linkedin.com/in/lorcan-mulvey-37ba7541"
$ws.Range("I6").Value = "lorcan.mulvey@mcaleer-rushe.co.uk , lorcanmulvey@yahoo.ie , lorcan.mulvey@yahoo.ie , lorcan.mulvey@berkeleygroup.co.uk"

$ws.Range("D7").Value = "https://linkedin.com/in/gerald-laverty-55a70790"
$ws.Range("I7").Value = ""

$ws.Range("D8").Value = "https://linkedin.com/in/darragh-greenan-5a2089b6"
$ws.Range("I8").Value = "darragh.greenan@mcaleer-rushe.co.uk"

$ws.Range("D9").Value = "https://linkedin.com/in/steve-morris-05a32933"
$ws.Range("I9").Value = ""

$ws.Range("D10").Value = "https://linkedin.com/in/lee-robert-gray-gradiosh-55b31a47"
$ws.Range("I10").Value = "leergray3@hotmail.co.uk , lee.gray@mcaleer-rushe.co.uk"

$ws.Range("D11").Value = "https://linkedin.com/in/eamonn-laverty-mciob-9634b9b1"
$ws.Range("I11").Value = ""

$ws.Range("D12").Value = "https://linkedin.com/in/eoin-gormley-10694a121"
$ws.Range("I12").Value = ""

$ws.Range("D13").Value = "https://linkedin.com/in/paddy-connolly-3a4527102"
$ws.Range("I13").Value = ""

$ws.Range("D14").Value = "https://linkedin.com/in/daisy-butterworth-14b437173"
$ws.Range("I14").Value = ""

$ws.Range("D15").Value = "https://linkedin.com/in/sinéad-gorman-she-her-523b9a64"
$ws.Range("I15").Value = ""

$ws.Range("D16").Value = "https://linkedin.com/in/connor-graham-50b2b241"
$ws.Range("I16").Value = "connor.graham@patton.co.uk , connor.graham@mcaleer-rushe.co.uk"

$ws.Range("D17").Value = "https://linkedin.com/in/cathal-magee-a397ba8a"
$ws.Range("I17").Value = "cathal.magee@mcaleer-rushe.co.uk , cathal.magee1@hotmail.co.uk"

$ws.Range("D18").Value = "https://linkedin.com/in/bruno-antoniazzi"
$ws.Range("I18").Value = "bruno.antoniazzi@mcaleer-rushe.co.uk"

$ws.Range("D19").Value = "https://linkedin.com/in/nina-salandy-bsc-hons-gradiosh-4675955a"
$ws.Range("I19").Value = ""

$ws.Range("D20").Value = "https://linkedin.com/in/peter-coyle-854a8063"
$ws.Range("I20").Value = ""

$ws.Range("D21").Value = "https://linkedin.com/in/orran-devine-b2715b166"
$ws.Range("I21").Value = ""

$ws.Range("D22").Value = "https://linkedin.com/in/niamh-heneghan-a98527197"
$ws.Range("I22").Value = ""
